# The workbook originally had a handful of name/place labels with a
# stray trailing space baked into the shared strings ("Ilaria ", "Botti ",
# "Di Marco ", "Disposal Station ", "Milano ", "Bologna "). Clean those up
# in place so the cell text matches the rest of the sheet (no trailing
# whitespace), then leave the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Ilaria"
$ws.Range("B7").Value = "Di Marco"
$ws.Range("A10").Value = "Disposal Station"
$ws.Range("A11").Value = "Disposal Station"
$ws.Range("B10").Value = "Milano"
$ws.Range("B11").Value = "Bologna"
$ws.Range("B5").Value = "Botti"

$ws.Range("A12").Select()
